$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '44.008.24'
$ws.Range('E2').Value = '  +0.61%  '
$ws.Range('D3').Value = '2.353.49'
$ws.Range('E3').Value = '  +0.09%  '
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.676'
$ws.Range('E5').Value = '  +0.50%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '239.06'
$ws.Range('E6').Value = '  +1.41%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '73.92'
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.582'
$ws.Range('E9').Value = '  +7.11%  '
$ws.Range('E10').Value = '  +1.15%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '57.19'
$ws.Range('E11').Value = '  +0.19%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '32.14'
$ws.Range('E12').Value = '  +13.52%  '
$ws.Range('E13').Value = '  +1.06%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.21'
$ws.Range('E14').Value = '  +7.51%  '
$ws.Range('D15').Value = '2.702.24'
$ws.Range('E15').Value = '  +0.10%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '16.57'
$ws.Range('E16').Value = '  -1.08%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.900'
$ws.Range('E17').Value = '  +1.04%  '
$ws.Range('D18').Value = '2.358.82'
$ws.Range('E18').Value = '  -2.35%  '
$ws.Range('D19').Value = '43.880.72'
$ws.Range('E19').Value = '  +0.29%  '
$ws.Range('E20').Value = '  +0.34%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.71'
$ws.Range('E21').Value = '  +4.26%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '76.80'
$ws.Range('E22').Value = '  -0.34%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '255.97'
$ws.Range('E23').Value = '  +0.96%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.92'
$ws.Range('E24').Value = '  +20.83%  '
$ws.Range('E25').Value = '  -0.05%  '
$ws.Range('E26').Value = '  -1.83%  '
$ws.Range('E27').Value = '  -0.32%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.67'
$ws.Range('E28').Value = '  +0.78%  '
$ws.Range('E29').Value = '  -1.93%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '22.69'
$ws.Range('E30').Value = '  +1.33%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '175.31'
$ws.Range('E31').Value = '  +1.48%  '
$ws.Range('E32').Value = '  -3.25%  '
$ws.Range('E33').Value = '  +2.67%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0756'
$ws.Range('E34').Value = '  +5.99%  '
$ws.Range('E35').Value = '  +1.98%  '
$ws.Range('E36').Value = '  +3.36%  '
$ws.Range('E37').Value = '  -3.38%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.36'
$ws.Range('E38').Value = '  -2.56%  '
$ws.Range('E39').Value = '  -1.25%  '
$ws.Range('E40').Value = '  +4.46%  '
$ws.Range('E41').Value = '  +12.31%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '19.21'
$ws.Range('E42').Value = '  -1.62%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.204'
$ws.Range('E43').Value = '  +12.38%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '9.08'
$ws.Range('E44').Value = '  +2.60%  '
$ws.Range('E45').Value = '  +0.12%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.68'
$ws.Range('E46').Value = '  +5.24%  '
$ws.Range('B47').Value = 'NEARProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.48'
$ws.Range('E47').Value = '  +8.41%  '
$ws.Range('B48').Value = 'TrustWalletToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.24'
$ws.Range('E48').Value = '  +1.33%  '
$ws.Range('B49').Value = 'MultiversX'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '57.34'
$ws.Range('E49').Value = '  +9.42%  '
$ws.Range('E50').Value = '  +0.50%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '99.94'
$ws.Range('E51').Value = '  +2.87%  '
